# Apply updated crypto price/volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.344.27"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "'1.842.32"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'240.04"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'0.6293"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'0.07413"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'24.86"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'0.07737"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "'1.829.09"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "'4.982"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'0.6794"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'0.00001024"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "'82.00"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "'6.254"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "'29.323.91"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'229.33"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'7.441"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'158.34"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'8.477"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").Value = "'17.44"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "'0.06518"
$ws.Range("E28").Value = "  +14.06%  "
$ws.Range("D29").Value = "'1.452"
$ws.Range("D30").Value = "'1.489"
$ws.Range("D31").Value = "'4.071"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'1.837"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").Value = "'0.6969"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "'0.01860"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "'2.816"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "'1.240.89"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").Value = "'6.796"
$ws.Range("E40").Value = "  +4.19%  "
$ws.Range("D41").Value = "'0.9343"
$ws.Range("D42").Value = "'0.9989"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "'100.86"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").Value = "'65.57"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "'7.056"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "'1.715"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.008"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1153"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'0.3904"
$ws.Range("E51").Value = "  -1.91%  "
